$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Job-Role column (D) with corrected / new role names
$ws.Range("D2").Value = "Sales Manager"
$ws.Range("D4").Value = "Business Development Manager"
$ws.Range("D6").Value = "Technical Manager"
$ws.Range("D5").Value = "Solution"

# Update the active cell selection to C9 (matches the saved view state)
$ws.Range("C9").Select()
